# Correction in SA algorithm and 746 logs
# Update the Fitness column (C) values:
#   - rows 2..101  (Generation 0..99)   -> 7573
#   - rows 102..140 (Generation 100..138) unchanged (already 7573)
#   - rows 141..166 (Generation 139..164) -> 7295
#   - rows 167..252 (Generation 165..250) -> 7293

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 101; $r++) {
    $ws.Cells.Item($r, 3).Value = 7573
}

for ($r = 141; $r -le 166; $r++) {
    $ws.Cells.Item($r, 3).Value = 7295
}

for ($r = 167; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
